$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (product_name), shifting
# product_name/account_name/quantity_on_hand/adjusted_quantity/description right.
$ws.Range("B1").EntireColumn.Insert()

# New "id" column: header + product id value, formatted like the other
# text columns (Calibri / General) and narrower to fit a short integer.
$ws.Columns("B:B").NumberFormat = "General"
$ws.Range("B1").Value = "id"
$ws.Range("B2").Value = 1827
$ws.Columns("B:B").ColumnWidth = 4.14

[void]$ws.Range("B3").Select()
